$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $style = $ws.Range($ref).Style
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).Style = $style
}

Set-TextValue $ws 'D2' '29.032.49'
Set-TextValue $ws 'E2' '  -0.45%  '
Set-TextValue $ws 'D3' '1.830.22'
Set-TextValue $ws 'E3' '  -0.32%  '
Set-TextValue $ws 'E4' '  +0.01%  '
Set-TextValue $ws 'D5' '241.43'
Set-TextValue $ws 'E5' '  -0.10%  '
Set-TextValue $ws 'D6' '0.6274'
Set-TextValue $ws 'E6' '  -5.03%  '
Set-TextValue $ws 'E7' '  +0.04%  '
Set-TextValue $ws 'D8' '0.07611'
Set-TextValue $ws 'E8' '  +2.38%  '
Set-TextValue $ws 'D9' '0.2920'
Set-TextValue $ws 'E9' '  -0.48%  '
Set-TextValue $ws 'D10' '22.82'
Set-TextValue $ws 'E10' '  -0.44%  '
Set-TextValue $ws 'D11' '0.07646'
Set-TextValue $ws 'E11' '  -1.38%  '
Set-TextValue $ws 'D12' '1.832.76'
Set-TextValue $ws 'E12' '  -0.51%  '
Set-TextValue $ws 'D13' '4.959'
Set-TextValue $ws 'E13' '  -0.87%  '
Set-TextValue $ws 'D14' '0.6656'
Set-TextValue $ws 'E14' '  -0.26%  '
Set-TextValue $ws 'D15' '82.45'
Set-TextValue $ws 'E15' '  -0.85%  '
Set-TextValue $ws 'D16' '0.000009471'
Set-TextValue $ws 'E16' '  +10.46%  '
Set-TextValue $ws 'D17' '5.988'
Set-TextValue $ws 'E17' '  -2.31%  '
Set-TextValue $ws 'D18' '28.951.99'
Set-TextValue $ws 'E18' '  -0.76%  '
Set-TextValue $ws 'D19' '225.00'
Set-TextValue $ws 'E19' '  -0.82%  '
Set-TextValue $ws 'D20' '12.33'
Set-TextValue $ws 'E20' '  -1.12%  '
Set-TextValue $ws 'D21' '0.9995'
Set-TextValue $ws 'E21' '  -0.19%  '
Set-TextValue $ws 'D22' '7.232'
Set-TextValue $ws 'E22' '  +1.73%  '
Set-TextValue $ws 'E23' '  +0.08%  '
Set-TextValue $ws 'D24' '161.08'
Set-TextValue $ws 'E24' '  +0.79%  '
Set-TextValue $ws 'D25' '8.422'
Set-TextValue $ws 'E25' '  -2.39%  '
Set-TextValue $ws 'E26' '  -2.60%  '
Set-TextValue $ws 'E27' '  -0.86%  '
Set-TextValue $ws 'D28' '1.493'
Set-TextValue $ws 'E28' '  -1.58%  '
Set-TextValue $ws 'D29' '4.054'
Set-TextValue $ws 'E29' '  -1.36%  '
Set-TextValue $ws 'D30' '4.033'
Set-TextValue $ws 'E30' '  -0.03%  '
Set-TextValue $ws 'D31' '1.199'
Set-TextValue $ws 'E31' '  +0.97%  '
Set-TextValue $ws 'D32' '0.05202'
Set-TextValue $ws 'E32' '  -1.63%  '
Set-TextValue $ws 'D33' '1.851'
Set-TextValue $ws 'E33' '  -0.84%  '
Set-TextValue $ws 'E34' '  +0.38%  '
Set-TextValue $ws 'D35' '0.7299'
Set-TextValue $ws 'E35' '  -1.25%  '
Set-TextValue $ws 'D36' '2.603'
Set-TextValue $ws 'E36' '  -2.22%  '
Set-TextValue $ws 'D37' '1.275.57'
Set-TextValue $ws 'E37' '  -1.86%  '
Set-TextValue $ws 'E38' '  +0.66%  '
Set-TextValue $ws 'D39' '0.01785'
Set-TextValue $ws 'E39' '  -0.39%  '
Set-TextValue $ws 'D40' '6.496'
Set-TextValue $ws 'E40' '  +7.46%  '
Set-TextValue $ws 'D41' '0.8926'
Set-TextValue $ws 'E41' '  -3.01%  '
Set-TextValue $ws 'E42' '  +0.09%  '
Set-TextValue $ws 'D43' '101.51'
Set-TextValue $ws 'E43' '  -0.40%  '
Set-TextValue $ws 'D44' '1.975.31'
Set-TextValue $ws 'E44' '  -0.78%  '
Set-TextValue $ws 'E45' '  -0.47%  '
Set-TextValue $ws 'D46' '63.79'
Set-TextValue $ws 'E46' '  +0.07%  '
Set-TextValue $ws 'D47' '0.00000000120'
Set-TextValue $ws 'E47' '  -1.17%  '
Set-TextValue $ws 'D48' '0.3982'
Set-TextValue $ws 'E48' '  -0.77%  '
Set-TextValue $ws 'D49' '0.07338'
Set-TextValue $ws 'E49' '  -12.40%  '
Set-TextValue $ws 'D50' '8.846'
Set-TextValue $ws 'E50' '  +1.38%  '
Set-TextValue $ws 'B51' 'Cronos'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D51' '0.05760'
Set-TextValue $ws 'E51' '  -1.51%  '
